$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 7).Value = 13.35945066666667
$ws.Cells.Item(2, 8).Value = 40.078352
$ws.Cells.Item(2, 9).Value = 0.04266941919869369
$ws.Cells.Item(2, 10).Value = 0.0426694191986937
$ws.Cells.Item(2, 13).Value = 15.03663066666667
$ws.Cells.Item(2, 14).Value = 45.109892
$ws.Cells.Item(2, 15).Value = 0.279146411176606
$ws.Cells.Item(2, 16).Value = 0.279146411176606
$ws.Cells.Item(2, 17).Value = 200.8811255842205
$ws.Cells.Item(2, 18).Value = 1807.930130257984
$ws.Cells.Item(2, 19).Value = 0.01191101523630552
$ws.Cells.Item(2, 20).Value = 0.01191101523630552
# Row 3
$ws.Cells.Item(3, 7).Value = 13.35945066666667
$ws.Cells.Item(3, 8).Value = 40.078352
$ws.Cells.Item(3, 9).Value = 0.04266941919869369
$ws.Cells.Item(3, 10).Value = 0.0426694191986937
$ws.Cells.Item(3, 15).Value = 0.2673306493381863
$ws.Cells.Item(3, 16).Value = 0.2673306493381863
$ws.Cells.Item(3, 17).Value = 192.3781914869049
$ws.Cells.Item(3, 18).Value = 1731.403723382144
$ws.Cells.Item(3, 19).Value = 0.01140684354127006
$ws.Cells.Item(3, 20).Value = 0.01140684354127006
# Row 4
$ws.Cells.Item(4, 7).Value = 13.35945066666667
$ws.Cells.Item(4, 8).Value = 40.078352
$ws.Cells.Item(4, 9).Value = 0.04266941919869369
$ws.Cells.Item(4, 10).Value = 0.0426694191986937
$ws.Cells.Item(4, 13).Value = 22.16851266666667
$ws.Cells.Item(4, 14).Value = 66.505538
$ws.Cells.Item(4, 15).Value = 0.411545703901694
$ws.Cells.Item(4, 16).Value = 0.411545703901694
$ws.Cells.Item(4, 17).Value = 296.1591513237085
$ws.Cells.Item(4, 18).Value = 2665.432361913376
$ws.Cells.Item(4, 19).Value = 0.01756041615920285
$ws.Cells.Item(4, 20).Value = 0.01756041615920285
# Row 5
$ws.Cells.Item(5, 7).Value = 13.35945066666667
$ws.Cells.Item(5, 8).Value = 40.078352
$ws.Cells.Item(5, 9).Value = 0.04266941919869369
$ws.Cells.Item(5, 10).Value = 0.0426694191986937
$ws.Cells.Item(5, 13).Value = 2.261165333333333
$ws.Cells.Item(5, 14).Value = 6.783496
$ws.Cells.Item(5, 15).Value = 0.04197723558351375
$ws.Cells.Item(5, 16).Value = 0.04197723558351374
$ws.Cells.Item(5, 17).Value = 30.20792671984355
$ws.Cells.Item(5, 18).Value = 271.871340478592
$ws.Cells.Item(5, 19).Value = 0.00179114426191527
$ws.Cells.Item(5, 20).Value = 0.00179114426191527
# Row 6
$ws.Cells.Item(6, 9).Value = 0.7730397221570005
$ws.Cells.Item(6, 10).Value = 0.7730397221570006
$ws.Cells.Item(6, 13).Value = 15.03663066666667
$ws.Cells.Item(6, 14).Value = 45.109892
$ws.Cells.Item(6, 15).Value = 0.279146411176606
$ws.Cells.Item(6, 16).Value = 0.279146411176606
$ws.Cells.Item(6, 17).Value = 3639.353251683478
$ws.Cells.Item(6, 18).Value = 32754.1792651513
$ws.Cells.Item(6, 19).Value = 0.2157912641370874
$ws.Cells.Item(6, 20).Value = 0.2157912641370873
# Row 7
$ws.Cells.Item(7, 9).Value = 0.7730397221570005
$ws.Cells.Item(7, 10).Value = 0.7730397221570006
$ws.Cells.Item(7, 15).Value = 0.2673306493381863
$ws.Cells.Item(7, 16).Value = 0.2673306493381863
$ws.Cells.Item(7, 19).Value = 0.2066572108884421
$ws.Cells.Item(7, 20).Value = 0.2066572108884421
# Row 8
$ws.Cells.Item(8, 9).Value = 0.7730397221570005
$ws.Cells.Item(8, 10).Value = 0.7730397221570006
$ws.Cells.Item(8, 13).Value = 22.16851266666667
$ws.Cells.Item(8, 14).Value = 66.505538
$ws.Cells.Item(8, 15).Value = 0.411545703901694
$ws.Cells.Item(8, 16).Value = 0.411545703901694
$ws.Cells.Item(8, 17).Value = 5365.50045332095
$ws.Cells.Item(8, 18).Value = 48289.50407988855
$ws.Cells.Item(8, 19).Value = 0.3181411765990728
$ws.Cells.Item(8, 20).Value = 0.3181411765990728
# Row 9
$ws.Cells.Item(9, 9).Value = 0.7730397221570005
$ws.Cells.Item(9, 10).Value = 0.7730397221570006
$ws.Cells.Item(9, 13).Value = 2.261165333333333
$ws.Cells.Item(9, 14).Value = 6.783496
$ws.Cells.Item(9, 15).Value = 0.04197723558351375
$ws.Cells.Item(9, 16).Value = 0.04197723558351374
$ws.Cells.Item(9, 17).Value = 547.2754894953387
$ws.Cells.Item(9, 18).Value = 4925.479405458048
$ws.Cells.Item(9, 19).Value = 0.03245007053239842
$ws.Cells.Item(9, 20).Value = 0.03245007053239842
# Row 10
$ws.Cells.Item(10, 7).Value = 10.23810666666667
$ws.Cells.Item(10, 8).Value = 30.71432
$ws.Cells.Item(10, 9).Value = 0.03270000212291218
$ws.Cells.Item(10, 10).Value = 0.03270000212291219
$ws.Cells.Item(10, 13).Value = 15.03663066666667
$ws.Cells.Item(10, 14).Value = 45.109892
$ws.Cells.Item(10, 15).Value = 0.279146411176606
$ws.Cells.Item(10, 16).Value = 0.279146411176606
$ws.Cells.Item(10, 17).Value = 153.9466286726044
$ws.Cells.Item(10, 18).Value = 1385.51965805344
$ws.Cells.Item(10, 19).Value = 0.009128088238078334
$ws.Cells.Item(10, 20).Value = 0.009128088238078334
# Row 11
$ws.Cells.Item(11, 7).Value = 10.23810666666667
$ws.Cells.Item(11, 8).Value = 30.71432
$ws.Cells.Item(11, 9).Value = 0.03270000212291218
$ws.Cells.Item(11, 10).Value = 0.03270000212291219
$ws.Cells.Item(11, 15).Value = 0.2673306493381863
$ws.Cells.Item(11, 16).Value = 0.2673306493381863
$ws.Cells.Item(11, 17).Value = 147.4303467954489
$ws.Cells.Item(11, 18).Value = 1326.87312115904
$ws.Cells.Item(11, 19).Value = 0.008741712800878184
$ws.Cells.Item(11, 20).Value = 0.008741712800878185
# Row 12
$ws.Cells.Item(12, 7).Value = 10.23810666666667
$ws.Cells.Item(12, 8).Value = 30.71432
$ws.Cells.Item(12, 9).Value = 0.03270000212291218
$ws.Cells.Item(12, 10).Value = 0.03270000212291219
$ws.Cells.Item(12, 13).Value = 22.16851266666667
$ws.Cells.Item(12, 14).Value = 66.505538
$ws.Cells.Item(12, 15).Value = 0.411545703901694
$ws.Cells.Item(12, 16).Value = 0.411545703901694
$ws.Cells.Item(12, 17).Value = 226.9635973226845
$ws.Cells.Item(12, 18).Value = 2042.67237590416
$ws.Cells.Item(12, 19).Value = 0.01345754539126078
$ws.Cells.Item(12, 20).Value = 0.01345754539126078
# Row 13
$ws.Cells.Item(13, 7).Value = 10.23810666666667
$ws.Cells.Item(13, 8).Value = 30.71432
$ws.Cells.Item(13, 9).Value = 0.03270000212291218
$ws.Cells.Item(13, 10).Value = 0.03270000212291219
$ws.Cells.Item(13, 13).Value = 2.261165333333333
$ws.Cells.Item(13, 14).Value = 6.783496
$ws.Cells.Item(13, 15).Value = 0.04197723558351375
$ws.Cells.Item(13, 16).Value = 0.04197723558351374
$ws.Cells.Item(13, 17).Value = 23.15005187363555
$ws.Cells.Item(13, 18).Value = 208.35046686272
$ws.Cells.Item(13, 19).Value = 0.001372655692694884
$ws.Cells.Item(13, 20).Value = 0.001372655692694884
# Row 14
$ws.Cells.Item(14, 7).Value = 47.46187333333333
$ws.Cells.Item(14, 8).Value = 142.38562
$ws.Cells.Item(14, 9).Value = 0.1515908565213935
$ws.Cells.Item(14, 10).Value = 0.1515908565213935
$ws.Cells.Item(14, 13).Value = 15.03663066666667
$ws.Cells.Item(14, 14).Value = 45.109892
$ws.Cells.Item(14, 15).Value = 0.279146411176606
$ws.Cells.Item(14, 16).Value = 0.279146411176606
$ws.Cells.Item(14, 17).Value = 713.6666600614489
$ws.Cells.Item(14, 18).Value = 6422.99994055304
$ws.Cells.Item(14, 19).Value = 0.0423160435651348
$ws.Cells.Item(14, 20).Value = 0.0423160435651348
# Row 15
$ws.Cells.Item(15, 7).Value = 47.46187333333333
$ws.Cells.Item(15, 8).Value = 142.38562
$ws.Cells.Item(15, 9).Value = 0.1515908565213935
$ws.Cells.Item(15, 10).Value = 0.1515908565213935
$ws.Cells.Item(15, 15).Value = 0.2673306493381863
$ws.Cells.Item(15, 16).Value = 0.2673306493381863
$ws.Cells.Item(15, 17).Value = 683.4584433347378
$ws.Cells.Item(15, 18).Value = 6151.125990012641
$ws.Cells.Item(15, 19).Value = 0.04052488210759596
$ws.Cells.Item(15, 20).Value = 0.04052488210759596
# Row 16
$ws.Cells.Item(16, 7).Value = 47.46187333333333
$ws.Cells.Item(16, 8).Value = 142.38562
$ws.Cells.Item(16, 9).Value = 0.1515908565213935
$ws.Cells.Item(16, 10).Value = 0.1515908565213935
$ws.Cells.Item(16, 13).Value = 22.16851266666667
$ws.Cells.Item(16, 14).Value = 66.505538
$ws.Cells.Item(16, 15).Value = 0.411545703901694
$ws.Cells.Item(16, 16).Value = 0.411545703901694
$ws.Cells.Item(16, 17).Value = 1052.159140173729
$ws.Cells.Item(16, 18).Value = 9469.432261563559
$ws.Cells.Item(16, 19).Value = 0.0623865657521576
$ws.Cells.Item(16, 20).Value = 0.06238656575215759
# Row 17
$ws.Cells.Item(17, 7).Value = 47.46187333333333
$ws.Cells.Item(17, 8).Value = 142.38562
$ws.Cells.Item(17, 9).Value = 0.1515908565213935
$ws.Cells.Item(17, 10).Value = 0.1515908565213935
$ws.Cells.Item(17, 13).Value = 2.261165333333333
$ws.Cells.Item(17, 14).Value = 6.783496
$ws.Cells.Item(17, 15).Value = 0.04197723558351375
$ws.Cells.Item(17, 16).Value = 0.04197723558351374
$ws.Cells.Item(17, 17).Value = 107.3191426363911
$ws.Cells.Item(17, 18).Value = 965.8722837275199
$ws.Cells.Item(17, 19).Value = 0.006363365096505166
$ws.Cells.Item(17, 20).Value = 0.006363365096505165

$wb.Save()